$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 377, shifting existing rows 377:395 down to 378:396
$ws.Rows("377:377").Insert()

# Populate the new row 377 with the weekly price-report entry
$ws.Range("A377").Value = 9
$ws.Range("B377").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C377").Value = "Metropolitana"
$ws.Range("D377").Value = 44753
$ws.Range("E377").Value = 13
$ws.Range("F377").Value = 100112052
$ws.Range("G377").Value = "Albahaca"
$ws.Range("H377").Value = "Sin especificar"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 430
$ws.Range("K377").Value = 4000
$ws.Range("L377").Value = 4000
$ws.Range("M377").Value = 4000
$ws.Range("N377").Value = "`$/paquete"
$ws.Range("O377").Value = "Región de Arica y Parinacota"
$ws.Range("P377").Value = 4000
$ws.Range("Q377").Value = 1
$ws.Range("R377").Value = "Hortaliza"
